$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("B9").Value = "PUT"
$ws.Range("C9").Value = "/pedidos/recibir"
$ws.Range("D9").Value = "Recibe un pedido"

$ws.Range("D10").Select()
